$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bugfix: Grande-Anse was incorrectly excluded when parsing the 2005 tax base
# data. Row 39 previously (wrongly) repeated "Grand Manan" in column A; it
# should instead be the "Grande-Anse" municipality with corrected totals.
# Row 38 ("Grand Manan") also gets corrected totals now that Grande-Anse's
# figures are no longer folded into it.

# Row 38: Grand Manan (name unchanged, values corrected)
$ws.Cells.Item(38, 2).Value = 103165300
$ws.Cells.Item(38, 3).Value = 552000
$ws.Cells.Item(38, 4).Value = 5635600
$ws.Cells.Item(38, 5).Value = 109352900
$ws.Cells.Item(38, 6).Value = 15022300
$ws.Cells.Item(38, 7).Value = 866300
$ws.Cells.Item(38, 8).Value = 1690700
$ws.Cells.Item(38, 9).Value = 17579300
$ws.Cells.Item(38, 10).Value = 126932200
$ws.Cells.Item(38, 11).Value = 135721850
$ws.Cells.Item(38, 12).Value = 135682764

# Row 39: fix mislabeled municipality name and its values -> Grande-Anse
$ws.Cells.Item(39, 1).Value = "Grande-Anse"
$ws.Cells.Item(39, 2).Value = 20258600
$ws.Cells.Item(39, 3).Value = 0
$ws.Cells.Item(39, 4).Value = 3606700
$ws.Cells.Item(39, 5).Value = 23865300
$ws.Cells.Item(39, 6).Value = 4262700
$ws.Cells.Item(39, 7).Value = 82300
$ws.Cells.Item(39, 8).Value = 472300
$ws.Cells.Item(39, 9).Value = 4817300
$ws.Cells.Item(39, 10).Value = 28682600
$ws.Cells.Item(39, 11).Value = 31091250
$ws.Cells.Item(39, 12).Value = 31070375
